# Append 8 new departure rows (307-314) to the "Main Data" sheet, as part of
# the "simple graph added, requires small fixes. added todo tasks" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ A=306; B="Saturday, Jan 14"; C="5:40 PM"; D="FR3263"; E="Dublin";    F="(DUB)"; G="Ryanair ";      H="B738"; I="(SP-RKM)"; J="5:51 PM"; L="0 hours, 11 minutes" },
    @{ A=307; B="Saturday, Jan 14"; C="6:00 PM"; D="W61773"; E="Reykjavik"; F="(KEF)"; G="Wizz Air ";     H="A320"; I="(HA-LYH)"; J="6:09 PM"; L="0 hours, 9 minutes" },
    @{ A=308; B="Saturday, Jan 14"; C="6:05 PM"; D="W61627"; E="Edinburgh"; F="(EDI)"; G="Wizz Air ";     H="A320"; I="(HA-LYO)"; J="6:51 PM"; L="0 hours, 46 minutes" },
    @{ A=309; B="Saturday, Jan 14"; C="6:30 PM"; D="FR6082"; E="London";    F="(STN)"; G="Ryanair ";      H="B738"; I="(SP-RSO)"; J="6:49 PM"; L="0 hours, 19 minutes" },
    @{ A=310; B="Saturday, Jan 14"; C="6:30 PM"; D="W61751"; E="Turku";     F="(TKU)"; G="Wizz Air ";     H="A321"; I="(HA-LTB)"; J="6:30 PM"; L="0 hours, 0 minutes" },
    @{ A=311; B="Saturday, Jan 14"; C="6:50 PM"; D="FR3279"; E="Milan";     F="(BGY)"; G="Ryanair ";      H="B738"; I="(SP-RKQ)"; J="7:05 PM"; L="0 hours, 15 minutes" },
    @{ A=312; B="Saturday, Jan 14"; C="7:25 PM"; D="FR1662"; E="Vienna";    F="(VIE)"; G="Lauda Europe "; H="A320"; I="(9H-LOS)"; J="7:28 PM"; L="0 hours, 3 minutes" },
    @{ A=313; B="Saturday, Jan 14"; C="7:55 PM"; D="W61611"; E="Liverpool"; F="(LPL)"; G="Wizz Air ";     H="A321"; I="(HA-LXL)"; J="9:03 PM"; L="1 hours, 8 minutes" }
)

$startRow = 307
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("G$r").Value = $row.G
    $ws.Range("H$r").Value = $row.H
    $ws.Range("I$r").Value = $row.I
    $ws.Range("J$r").Value = $row.J
    # K column (DIFFERENCE label) stays blank for data rows, like the rest of the table.
    $ws.Range("K$r").Borders.LineStyle = 0
    $ws.Range("L$r").Value = $row.L
    # M column is always blank in this table.
    $ws.Range("M$r").Borders.LineStyle = 0
}
